$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update story text (column A), locationEvent (column B), socialEvent (column C)
# and row heights for rows 22-33 to reflect the rewritten proposal scene. ---

$ws.Range("A22").Value = "The professor stopped the lecture and asked Lana and Aaron to come to the front of the class to do a demonstration. "
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = 2
$ws.Rows.Item(22).RowHeight = 92.05

$ws.Range("A23").Value = "They walked to the front of the class with Lana very confused, and when they got to the front Aaron got down on one knee and asked Lana to marry him. "
$ws.Range("B23").Value = 3
$ws.Range("C23").Value = 3
$ws.Rows.Item(23).RowHeight = 116.05

$ws.Range("A24").Value = "The professor quickly switched to the last slide on the powerpoint which was a photo collage Aaron had sent him of pictures of Aaron and Lana. "
$ws.Range("B24").Value = 3
$ws.Range("C24").Value = 3
$ws.Rows.Item(24).RowHeight = 116.05

$ws.Range("A25").Value = "He opened the small notebook he had in his hand and the ring was sitting on the page. "
$ws.Range("B25").Value = 3
$ws.Range("C25").Value = 3
$ws.Rows.Item(25).RowHeight = 80.05

$ws.Range("A26").Value = "Lana had put her hand over her mouth and tears were streaming down her face, but finally, she said, " + [char]8220 + "Yes!" + [char]8221
$ws.Range("B26").Value = 3
$ws.Range("C26").Value = 4
$ws.Rows.Item(26).RowHeight = 92.05

$ws.Range("A27").Value = "Aaron put the ring on her finger, and they hugged while the whole class applauded and cheered. "
$ws.Range("B27").Value = 3
$ws.Range("C27").Value = 4
$ws.Rows.Item(27).RowHeight = 80.05

$ws.Range("A28").Value = "Their friends, who had been secretly waiting outside, came into the class and embraced both of them."
$ws.Range("B28").Value = 3
$ws.Range("C28").Value = 4
$ws.Rows.Item(28).RowHeight = 92.05

$ws.Range("A29").Value = "The professor gave Aaron a celebratory handshake."
$ws.Range("B29").Value = 3
$ws.Range("C29").Value = 4
$ws.Rows.Item(29).RowHeight = 56.05

$ws.Range("A30").Value = "They made their way back to their seats, and Lana pulled out her phone to text her sister."
$ws.Range("B30").Value = 3
$ws.Range("C30").Value = 4
$ws.Rows.Item(30).RowHeight = 80.05

$ws.Range("A31").Value = "The professor dismissed the class, saying, " + [char]8220 + "That was enough excitement for today!"
$ws.Range("B31").Value = 4
$ws.Range("C31").Value = 4
$ws.Rows.Item(31).RowHeight = 80.05

$ws.Range("A32").Value = "Don" + [char]8217 + "t forget there" + [char]8217 + "s a midterm exam next week on Monday on the material covered in the first 4 weeks of the semester."
$ws.Range("B32").Value = 4
$ws.Range("C32").Value = 4
$ws.Rows.Item(32).RowHeight = 104.05

$ws.Range("A33").Value = "As everyone was packing up their notebooks and leaving, Lana and Aaron were surrounded by their friends eager to get a glimpse of the ring and congratulate them. "
$ws.Range("B33").Value = 4
$ws.Range("C33").Value = 4
$ws.Rows.Item(33).RowHeight = 140.05

# --- The old rows 34-36 (now redundant story text, since the demonstration/proposal
# got condensed into fewer rows above) are removed. Deleting them shifts the trailing
# blank padding rows (previously 37-45) up into place, which also naturally drops the
# sheet from 45 rows down to 42 rows. ---

$ws.Rows.Item(34).Delete()
$ws.Rows.Item(34).Delete()
$ws.Rows.Item(34).Delete()
